# Grid Hack.docx - apply proofing-pass XML restructure + To-Do list reorder
# described by the target diff. Since the Word object model here has no
# direct "insert w:proofErr" command, runs are rebuilt verbatim (including
# w:proofErr spellStart/spellEnd/gramStart/gramEnd markers) via
# Range.InsertXML, which accepts a WordprocessingML package fragment and
# replaces the target range's content with it.

function Set-RangeXml {
    param($range, [string]$innerXml)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Set-ParagraphsXml {
    param($doc, [int]$firstIndex, [int]$lastIndex, [string]$innerXml)
    $startRange = $doc.Paragraphs($firstIndex).Range
    $endRange = $doc.Paragraphs($lastIndex).Range
    $rng = $doc.Range($startRange.Start, $endRange.End)
    Set-RangeXml $rng $innerXml
}

$d = $word.ActiveDocument

# Work from the bottom of the document upward so that edits which change
# the paragraph count (the To-Do reorder below) never invalidate the
# paragraph indices of not-yet-processed text further up the document.

# --- Doors.js : "onMouseOver =  function(args)" -> wrap identifiers in proofErr ---
Set-ParagraphsXml $d 84 84 @'
<w:p><w:pPr><w:pStyle w:val="LTUHeading4"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>onMouseOver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> =  function</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

# --- Doors.js : "onUse = function(args)" -> wrap identifiers in proofErr ---
Set-ParagraphsXml $d 83 83 @'
<w:p><w:pPr><w:pStyle w:val="LTUHeading4"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>onUse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> = function</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

# --- Doors.js : "onUse" heading + "The args given to Door ..." paragraph ---
Set-ParagraphsXml $d 79 80 @'
<w:p><w:pPr><w:pStyle w:val="LTUHeading4"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>onUse</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> given to Door may require knowledge of the actual user.</w:t></w:r></w:p>
'@

# --- stateSetup.js : "setMouseFocus = function(cellX, cellY, mouseX, mouseY)" ---
Set-ParagraphsXml $d 63 63 @'
<w:p><w:pPr><w:pStyle w:val="LTUHeading4"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>setMouseFocus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> = function(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>cellX</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>cellY</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>mouseX</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>mouseY</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

# --- stateSetup.js : "addPlayer = function(args)" ---
Set-ParagraphsXml $d 62 62 @'
<w:p><w:pPr><w:pStyle w:val="LTUHeading4"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>addPlayer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> = function(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

# --- stateSetup.js : "msgPump" -> wrap in proofErr ---
Set-ParagraphsXml $d 56 56 @'
<w:p><w:pPr><w:pStyle w:val="LTUHeading4"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>msgPump</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

# --- stateSetup.js To Do: drop the two struck-through items, keep the
#     remaining three (moving the _GoBack bookmark onto the new first one) ---
Set-ParagraphsXml $d 48 52 @'
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Setup turn order prior to starting the game</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t>Set visible monsters to active</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t>Set invisible monsters to inactive</w:t></w:r></w:p>
'@

# --- gridHack To Do: "Consider an inventory interface in the hud ..." ---
Set-ParagraphsXml $d 44 44 @'
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Consider an inventory interface in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>hud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> for selecting weapons.  Maybe simply a form of drop down box.</w:t></w:r></w:p>
'@

# --- gridHack To Do: "The hud should not scale with the map ..." ---
Set-ParagraphsXml $d 39 39 @'
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>hud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> should not scale with the map: it has a fixed location.</w:t></w:r></w:p>
'@

# --- gridHack To Do: "Add a hud to the display ..." ---
Set-ParagraphsXml $d 38 38 @'
<w:p><w:pPr><w:pStyle w:val="LTUNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add a </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>hud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to the display to allow the user to define actions.</w:t></w:r></w:p>
'@

Write-Output "done"
